# "Fruta / hortaliza, semanal" — insert a new weekly price-survey record
# for Espinaca (Vega Modelo de Temuco) as row 37, pushing the existing
# rows 37-90 down to 38-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37; this shifts rows 37:90 -> 38:91
# and grows the sheet dimension from A1:R90 to A1:R91 automatically.
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the new survey record.
$ws.Cells.Item(37, 1).Value = 10
$ws.Cells.Item(37, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37, 3).Value = "La Araucanía"
$ws.Cells.Item(37, 4).Value = 44495
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(37, 6).Value = 100112012
$ws.Cells.Item(37, 7).Value = "Espinaca"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 20
$ws.Cells.Item(37, 11).Value = 8000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 8000
$ws.Cells.Item(37, 14).Value = "$/docena de atados"
$ws.Cells.Item(37, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(37, 16).Value = 2667
$ws.Cells.Item(37, 17).Value = 3
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Match the date-formatted number format used by the rest of column D.
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
